$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (rows 100-101),
# pushing the existing rows 100-202 down to 102-204.
$ws.Rows("100:101").Insert()

# Row 100: new weekly entry - Red Globe / Especial
$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 45033
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100109
$ws.Range("H100").Value = "Uva"
$ws.Range("I100").Value = 100109001
$ws.Range("J100").Value = "Uva"
$ws.Range("K100").Value = "Red Globe"
$ws.Range("L100").Value = "Especial"
$ws.Range("M100").Value = 60
$ws.Range("N100").Value = 12000
$ws.Range("O100").Value = 12000
$ws.Range("P100").Value = 12000
$ws.Range("Q100").Value = "$/bandeja 18 kilos"
$ws.Range("R100").Value = "Región de O'Higgins"
$ws.Range("S100").Value = 667
$ws.Range("T100").Value = 18

# Row 101: new weekly entry - Red Globe / Primera
$ws.Range("A101").Value = 7
$ws.Range("B101").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C101").Value = "Ñuble"
$ws.Range("D101").Value = 45033
$ws.Range("E101").Value = 16
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100109
$ws.Range("H101").Value = "Uva"
$ws.Range("I101").Value = 100109001
$ws.Range("J101").Value = "Uva"
$ws.Range("K101").Value = "Red Globe"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 60
$ws.Range("N101").Value = 10000
$ws.Range("O101").Value = 10000
$ws.Range("P101").Value = 10000
$ws.Range("Q101").Value = "$/bandeja 18 kilos"
$ws.Range("R101").Value = "Región de O'Higgins"
$ws.Range("S101").Value = 556
$ws.Range("T101").Value = 18
